$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B6: "Receive commands from surface" -> "Receive commands from surface, send commands to Arduino"
$ws.Range("B6").Value = "Receive commands from surface, send commands to Arduino"

# Fill in B7 (was empty) with new text "Reading sensor data"
$ws.Range("B7").Value = "Reading sensor data"

# Update the active cell selection from B9 to C9
$ws.Range("C9").Select()
